# Added date email sent to the top of each column, attachment row data
# split after ":" - read in as double if number, else read as string.
#
# Column A keeps the original "email #1" attachment rows (minus the
# "#22222222222222" line, which belongs to the second email) with a new
# date at the top (row 1). Column B holds the original column A content
# shifted down one row (including the "#22222222222222" line) under its
# own date header, with the former row-7 value now rounded to 2 dp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new column B width (Excel auto-fit style width for the new data) ---
$ws.Columns.Item(2).ColumnWidth = 15.75

# --- row 1: the two "date email sent" values -----------------------------
$ws.Range("A1").Value = 42289.739444444444
$ws.Range("A1").NumberFormat = "m/d/yy h:mm"

$ws.Range("B1").Value = 42299.564837962964
$ws.Range("B1").NumberFormat = "m/d/yy h:mm"

# --- column A: only the attachment-number row actually changes -----------
# (A2 "This is a heading.", A4 "More details...", A6 204.33 and A8 701.9
#  already hold the right values from the original file, so they are left
#  untouched; only A1 -- replaced by the date above -- and A7, now rounded
#  to 2dp, need writing.)
$ws.Range("A7").Value = 201.23

# --- column B: second email's attachment rows -----------------------------
$ws.Range("B2").Value = "#22222222222222"
$ws.Range("B3").Value = "This is a heading. "
$ws.Range("B5").Value = "More details. Test # 1. "
$ws.Range("B7").Value = 204.33
$ws.Range("B8").Value = 201.23
$ws.Range("B9").Value = 701.9

# --- selection / active cell housekeeping (matches Excel's post-edit state)
$ws.Range("A1:C10").Select() | Out-Null
